$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45879
$ws.Range("B2").Value = 104.72
$ws.Range("C2").Value = 85.44
$ws.Range("D2").Value = 79
$ws.Range("E2").Value = 68
$ws.Range("F2").Value = 85
$ws.Range("G2").Value = 90.5
$ws.Range("H2").Value = 96.44
$ws.Range("I2").Value = 90
$ws.Range("J2").Value = 49.9
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = -0.25
$ws.Range("M2").Value = -0.25
$ws.Range("N2").Value = -0.5
$ws.Range("O2").Value = -0.62
$ws.Range("P2").Value = -0.6
$ws.Range("Q2").Value = -0.6
$ws.Range("R2").Value = -0.5
$ws.Range("S2").Value = 1.72
$ws.Range("T2").Value = 26.5
$ws.Range("U2").Value = 80.5
$ws.Range("V2").Value = 98.01000000000001
$ws.Range("W2").Value = 130
$ws.Range("X2").Value = 136.87
$ws.Range("Y2").Value = 112.98
$ws.Range("Z2").Value = 55.51
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 119.46
$ws.Range("AC2").Value = "22h-24h"
$ws.Range("AD2").Value = 124.93
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 114
$ws.Range("AG2").Value = "8h-18h"
